$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark from the last (empty) paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert the new "_GoBack" bookmark at the very start of the document
# (right before the first run of the first paragraph), mirroring what
# Word does when the most recent edit happens at the top of the doc.
$start = $d.Range(0, 0)
$d.Bookmarks.Add("_GoBack", $start) | Out-Null
